$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Insert 6 new blank rows right after row 6 (i.e. at row 7), pushing the
# old row 7 ("CPU Utilization" ...) and everything below it down by 6.
# ---------------------------------------------------------------------
$ws.Rows("7:12").Insert()

# ---------------------------------------------------------------------
# New row 7: plain text, no special formatting.
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Pod CPU/Memory requests define a set amount of CPU and memory that the pod needs on a regular basis"

# ---------------------------------------------------------------------
# New row 8: merged A8:C8, same "note box" look used elsewhere in the
# sheet (left/top aligned, wrapped text) - copy format from a cell that
# already has that look (A13, old A19 analogue) then set text + height.
# ---------------------------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("A8:C8").Merge()
$ws.Range("A8").Value = "Pod CPU/Memory limits are the maximum amount of CPU and memory that a pod can use. These limits help define which pods should be killed in the event of node instability due to insufficient resources."
$ws.Rows("8:8").RowHeight = 30.5

# ---------------------------------------------------------------------
# New row 10: merged A10:C10 and D10:H10, a highlighted red-text note
# box with a border around the text (new style: red font + thin border
# + left/top aligned wrapped text).
# ---------------------------------------------------------------------
$ws.Range("A10:C10").Merge()
$ws.Range("D10:H10").Merge()
$ws.Range("A3").Copy()
$ws.Range("A10:H10").PasteSpecial(-4122)
$ws.Range("A10:H10").HorizontalAlignment = -4131
$ws.Range("A10:H10").VerticalAlignment = -4160
$ws.Range("A10:H10").Font.Color = 255
$ws.Range("A10").Value = "it is very important to monitor the performance of your application at different times during the day or week. Determine when the peak demand is, and align the pod limits to the resources required to meet the application's max needs."
$ws.Range("D10").Value = "For this activity to take place we have to deploy and test the entire functioning application so that the deployment files can be changed accordingly for requests and limits."
$ws.Rows("10:10").RowHeight = 44.5

# ---------------------------------------------------------------------
# Column D is slightly wider in the new layout.
# ---------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 13.08984375

# ---------------------------------------------------------------------
# Selection moves to the new A10:C10 box.
# ---------------------------------------------------------------------
$ws.Range("A10:C10").Select()
